$d = $word.ActiveDocument

# paragraph 92
$p = $d.Paragraphs.Item(92)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# paragraph 91
$p = $d.Paragraphs.Item(91)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D971ED" w:rsidRDefault="00CC2A52" w:rsidP="00CC2A52"><w:r w:rsidRPr="00CC2A52"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>到排队的时候提醒</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Eating time</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# paragraph 90
$p = $d.Paragraphs.Item(90)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00CC2A52" w:rsidRDefault="00CC2A52" w:rsidP="00CC2A52"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="00CC2A52"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>还有两队的时候提醒</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Eating time coming</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# delete paragraph 89 (取消排队 / Sequence cancled push message)
$p89 = $d.Paragraphs.Item(89)
$p89.Range.Delete()

# paragraph 88
$p = $d.Paragraphs.Item(88)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00CC2A52" w:rsidRDefault="00CC2A52" w:rsidP="00CC2A52"><w:r w:rsidRPr="00CC2A52"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>订单确认</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Order confirmed</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# paragraph 87
$p = $d.Paragraphs.Item(87)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00CC2A52" w:rsidRDefault="00CC2A52" w:rsidP="00CC2A52"><w:r w:rsidRPr="00CC2A52"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>饭局邀请反馈</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Invitation feedback</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# paragraph 86
$p = $d.Paragraphs.Item(86)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00CC2A52" w:rsidRDefault="00CC2A52" w:rsidP="00CC2A52"><w:r w:rsidRPr="00CC2A52"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>饭局邀请</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Meal Invitation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# paragraph 85
$p = $d.Paragraphs.Item(85)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00CC2A52" w:rsidRPr="00CC2A52" w:rsidRDefault="00CC2A52" w:rsidP="00CC2A52"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:ind w:firstLineChars="0"/><w:rPr><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00CC2A52"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>消息提醒</w:t></w:r><w:r w:rsidRPr="00CC2A52"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>Title</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# paragraph 82
$p = $d.Paragraphs.Item(82)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D971ED" w:rsidRPr="00D971ED" w:rsidRDefault="00D971ED" w:rsidP="00D971ED"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/><w:rPr><w:b/><w:color w:val="00B050"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00D971ED"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:color w:val="00B050"/><w:sz w:val="28"/></w:rPr><w:t>推送相关</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# paragraph 80
$p = $d.Paragraphs.Item(80)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0057140B" w:rsidRPr="00F373AD" w:rsidRDefault="00CC2A52" w:rsidP="00F373AD"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:ind w:firstLineChars="0"/><w:rPr><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00F373AD"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>获取排队详情</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

